# Update the NATMI ligand-receptor edge table with the refreshed TPM-based
# results. The sending/target cluster assignments and every numeric metric
# (columns E..T) for rows 2-7 are replaced with the new values; the table
# grows from 4 data rows (2-5) to 6 data rows (2-7).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each inner array is one data row (columns A..T) for sheet rows 2..7.
$data = @(
  @("ECs", "Nppc", "Npr2", "ECs", 1, 0.3333333333333333, 0.224617, 0.673851, 0.7269197784238318, 0.7269197784238318, 3, 1, 13.18988666666667, 39.56966, 0.6665003854349094, 0.6665003854349095, 2.962672773406666, 26.66405496066, 0.4844923124997428, 0.4844923124997429),
  @("ECs", "Nppc", "Npr2", "FAPs", 1, 0.3333333333333333, 0.224617, 0.673851, 0.7269197784238318, 0.7269197784238318, 3, 1, 5.523445333333334, 16.570336, 0.2791061467494529, 0.2791061467494529, 1.240659720437333, 11.165937483936, 0.2028877783518418, 0.2028877783518418),
  @("ECs", "Nppc", "Npr2", "MuSCs", 1, 0.3333333333333333, 0.224617, 0.673851, 0.7269197784238318, 0.7269197784238318, 3, 1, 1.076434, 3.229302, 0.05439346781563764, 0.05439346781563765, 0.241785375778, 2.176068382002, 0.03953968757224714, 0.03953968757224715),
  @("MuSCs", "Nppc", "Npr2", "ECs", 2, 0.6666666666666666, 0.08438133333333332, 0.253144, 0.2730802215761681, 0.2730802215761681, 3, 1, 13.18988666666667, 39.56966, 0.6665003854349094, 0.6665003854349095, 1.112980223448889, 10.01682201104, 0.1820080729351665, 0.1820080729351665),
  @("MuSCs", "Nppc", "Npr2", "FAPs", 2, 0.6666666666666666, 0.08438133333333332, 0.253144, 0.2730802215761681, 0.2730802215761681, 3, 1, 5.523445333333334, 16.570336, 0.2791061467494529, 0.2791061467494529, 0.4660756818204445, 4.194681136384, 0.07621836839761111, 0.07621836839761111),
  @("MuSCs", "Nppc", "Npr2", "MuSCs", 2, 0.6666666666666666, 0.08438133333333332, 0.253144, 0.2730802215761681, 0.2730802215761681, 3, 1, 1.076434, 3.229302, 0.05439346781563764, 0.05439346781563765, 0.09083093616533333, 0.817478425488, 0.0148537802433905, 0.0148537802433905),
)

$r = 2
foreach ($row in $data) {
  for ($c = 1; $c -le $row.Length; $c++) {
    $ws.Cells.Item($r, $c).Value = $row[$c - 1]
  }
  $r++
}
